$wb = $excel.ActiveWorkbook

# ALC row 26
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 6250
$ws.Range("I26").Value = 4000
$ws.Range("J26").Value = 7000
$ws.Range("K26").Value = 4000
$ws.Range("L26").Value = 7000
$ws.Range("M26").Value = -3656
$ws.Range("N26").Value = -7688

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I62").Value = 10000
$ws.Range("K62").Value = 10000
$ws.Range("M62").Value = -9376

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I65").Value = 10000
$ws.Range("K65").Value = 50000
$ws.Range("M65").Value = -46880

# ALC row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 1693.9048
$ws.Range("I103").Value = 2838.3333
$ws.Range("J103").Value = 1236.1333
$ws.Range("K103").Value = 8514.999899999999
$ws.Range("L103").Value = 3708.3999
$ws.Range("M103").Value = -7928.999899999999
$ws.Range("N103").Value = -4880.3999

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 6244.909
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2611.8708
$ws.Range("I137").Value = 1838.2
$ws.Range("K137").Value = 5514.6
$ws.Range("M137").Value = -2964.6

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 867
$ws.Range("I2").Value = 770.5833
$ws.Range("K2").Value = 770.5833
$ws.Range("M2").Value = -657.5833

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3033730
$ws.Range("I32").Value = 3449.923
$ws.Range("J32").Value = 14289056
$ws.Range("K32").Value = 3449.923
$ws.Range("L32").Value = 14289056
$ws.Range("M32").Value = -3162.923
$ws.Range("N32").Value = -14289630

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4358.25
$ws.Range("I45").Value = 3999.5
$ws.Range("K45").Value = 3999.5
$ws.Range("M45").Value = -3622.5

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 884
$ws.Range("I110").Value = 884
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 884
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1161
$ws.Range("N110").ClearContents()

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 867
$ws.Range("I116").Value = 770.5833
$ws.Range("K116").Value = 770.5833
$ws.Range("M116").Value = 1523.4167

# ARM row 118
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 867
$ws.Range("I3").Value = 770.5833
$ws.Range("K3").Value = 770.5833
$ws.Range("M3").Value = -656.5833

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1084.3334
$ws.Range("I20").Value = 969.9091
$ws.Range("J20").Value = 1399
$ws.Range("K20").Value = 969.9091
$ws.Range("L20").Value = 1399
$ws.Range("M20").Value = -722.9091
$ws.Range("N20").Value = -1893

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 5187.3335
$ws.Range("I107").Value = 2972.5
$ws.Range("K107").Value = 2972.5
$ws.Range("M107").Value = -1052.5

# BSM row 131
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4711.0376
$ws.Range("I31").Value = 2757.2856
$ws.Range("K31").Value = 2757.2856
$ws.Range("M31").Value = -2462.2856

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4711.0376
$ws.Range("I34").Value = 2757.2856
$ws.Range("K34").Value = 2757.2856
$ws.Range("M34").Value = -2555.2856

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2474.75
$ws.Range("I99").Value = 2299.5
$ws.Range("K99").Value = 2299.5
$ws.Range("M99").Value = -801.5

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 761.1667
$ws.Range("I107").Value = 625
$ws.Range("J107").Value = 829.25
$ws.Range("K107").Value = 625
$ws.Range("L107").Value = 829.25
$ws.Range("M107").Value = 1295
$ws.Range("N107").Value = -4669.25

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2474.75
$ws.Range("I126").Value = 2299.5
$ws.Range("K126").Value = 6898.5
$ws.Range("M126").Value = -4428.5

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 485897.6
$ws.Range("I4").Value = 800070.1
$ws.Range("J4").Value = 358.27274
$ws.Range("K4").Value = 2400210.3
$ws.Range("L4").Value = 1074.81822
$ws.Range("M4").Value = -2400098.3
$ws.Range("N4").Value = -1298.81822

# CUL row 7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 274.31818
$ws.Range("I7").Value = 140.76923
$ws.Range("J7").Value = 467.22223
$ws.Range("K7").Value = 422.30769
$ws.Range("L7").Value = 1401.66669
$ws.Range("M7").Value = -310.30769
$ws.Range("N7").Value = -1625.66669

# CUL row 9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1021
$ws.Range("I9").Value = 85
$ws.Range("J9").Value = 1333
$ws.Range("K9").Value = 255
$ws.Range("L9").Value = 3999
$ws.Range("M9").Value = -31
$ws.Range("N9").Value = -4447

# CUL row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 404.94116
$ws.Range("I38").Value = 399
$ws.Range("J38").Value = 500
$ws.Range("K38").Value = 1197
$ws.Range("L38").Value = 1500
$ws.Range("M38").Value = -850
$ws.Range("N38").Value = -2194

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2062.76
$ws.Range("I131").Value = 1032.6666
$ws.Range("K131").Value = 3097.9998
$ws.Range("M131").Value = 1942.0002

# CUL row 141
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 3920
$ws.Range("I141").Value = 1200
$ws.Range("K141").Value = 3600
$ws.Range("M141").Value = 1580

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2378.6316
$ws.Range("I107").Value = 1270.7
$ws.Range("K107").Value = 1270.7
$ws.Range("M107").Value = 649.3

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4780.154
$ws.Range("I113").Value = 2471.3333
$ws.Range("K113").Value = 2471.3333
$ws.Range("M113").Value = -301.3332999999998

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3087.9092
$ws.Range("I126").Value = 1980
$ws.Range("J126").Value = 4011.1667
$ws.Range("K126").Value = 5940
$ws.Range("L126").Value = 12033.5001
$ws.Range("M126").Value = -3470
$ws.Range("N126").Value = -16973.5001

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 937.3333
$ws.Range("I132").Value = 937.3333
$ws.Range("K132").Value = 2811.9999
$ws.Range("M132").Value = -281.9998999999998

# LTW row 3
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 4600
$ws.Range("I3").Value = 4750
$ws.Range("K3").Value = 4750
$ws.Range("M3").Value = -4638

# LTW row 15
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H15").Value = 4600
$ws.Range("I15").Value = 4750
$ws.Range("K15").Value = 4750
$ws.Range("M15").Value = -4580

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5507.4287
$ws.Range("I61").Value = 2999
$ws.Range("J61").Value = 7388.75
$ws.Range("K61").Value = 2999
$ws.Range("L61").Value = 7388.75
$ws.Range("M61").Value = -2797
$ws.Range("N61").Value = -7792.75

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 5507.4287
$ws.Range("I113").Value = 2999
$ws.Range("J113").Value = 7388.75
$ws.Range("K113").Value = 2999
$ws.Range("L113").Value = 7388.75
$ws.Range("M113").Value = -829
$ws.Range("N113").Value = -11728.75

# WVR row 11
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1550
$ws.Range("I100").Value = 1437.5
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 2875
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -2334
$ws.Range("N100").Value = -5082
